$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.167.90'
$ws.Range('E2').Value = '  +0.26%  '
$ws.Range('D3').Value = '2.547.61'
$ws.Range('E3').Value = '  +3.35%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.47'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.26'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +3.30%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.587'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.16%  '
$ws.Range('D9').Value = '2.543.34'
$ws.Range('E9').Value = '  +3.22%  '
$ws.Range('E10').Value = '  -0.02%  '
$ws.Range('E11').Value = '  -1.72%  '
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.352'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '27.51'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.82%  '
$ws.Range('D15').Value = '3.001.59'
$ws.Range('E15').Value = '  +3.32%  '
$ws.Range('D16').Value = '63.105.11'
$ws.Range('E16').Value = '  +0.34%  '
$ws.Range('E17').Value = '  +2.03%  '
$ws.Range('D18').Value = '2.543.06'
$ws.Range('E18').Value = '  +3.14%  '
$ws.Range('E19').Value = '  +2.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '335.85'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.29%  '
$ws.Range('E21').Value = '  +1.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.77'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E23').Value = '  +0.15%  '
$ws.Range('E24').Value = '  -0.33%  '
$ws.Range('B25').Value = 'Kaspa'
$ws.Range('C25').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.170'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.12%  '
$ws.Range('B26').Value = 'Fetch.AI'
$ws.Range('C26').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.63'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +8.64%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.52'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +10.07%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.53'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.93%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +7.78%  '
$ws.Range('D31').Value = '0.0₃0823'
$ws.Range('E31').Value = '  +2.74%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.84'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '175.82'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E34').Value = '  +4.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '411.73'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.54%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.400'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.89%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.40'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.31%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.15%  '
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '39.28'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.91%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '153.31'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.44%  '
$ws.Range('E44').Value = '  +2.03%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '20.95'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.93%  '
$ws.Range('E46').Value = '  +1.37%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0963'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0522'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.46%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.36'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.49%  '
$ws.Range('E51').Value = '  +2.80%  '
